$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three summary rows: 100 -> 0M, 0 -> 0M, 273 -> 0M
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# 2) Insert 10 new rows right before the row that used to be row 4 ("0"),
#    in order, by always adding a fresh row immediately above that same
#    anchor row (its index grows by one after each insert).
$newValues = @("136", "0.00002", "0.00005", "0.00004", "0.00001", "0.00003", "0.00004", "0.00005", "0.00471", "100.0")
$anchorIndex = 4
foreach ($val in $newValues) {
    $anchorRow = $t.Rows.Item($anchorIndex)
    $t.Rows.Add($anchorRow) | Out-Null
    $t.Cell($anchorIndex, 1).Range.Text = $val
    $anchorIndex = $anchorIndex + 1
}

# 3) The final three rows held 10 tab-separated values crammed into one
#    run; collapse each back down to just its first value.
$t.Cell(44,1).Range.Text = "100"
$t.Cell(45,1).Range.Text = "0"
$t.Cell(46,1).Range.Text = "273"
